# Apply the "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# Sheets (by index, matching workbook order):
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Performances)
#   3 = 本地生活 (Local Life)
#   4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: 展览 - update "want to go" counts (and one date range fix)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F3").Value  = 1825
$ws1.Range("F4").Value  = 405
$ws1.Range("F5").Value  = 1491
$ws1.Range("F6").Value  = 855
$ws1.Range("F8").Value  = 740
$ws1.Range("E9").Value  = "2024.06.08 10:00-06.09 16:00"
$ws1.Range("F9").Value  = 13187
$ws1.Range("F10").Value = 13063
$ws1.Range("F11").Value = 995
$ws1.Range("F15").Value = 67
$ws1.Range("F16").Value = 639
$ws1.Range("F17").Value = 2063
$ws1.Range("F18").Value = 54
$ws1.Range("F21").Value = 18
$ws1.Range("F22").Value = 193
$ws1.Range("F23").Value = 274
$ws1.Range("F24").Value = 736

# -----------------------------------------------------------------
# Sheet 2: 演出 - update "want to go" counts
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F4").Value = 15
$ws2.Range("F9").Value = 1

# -----------------------------------------------------------------
# Sheet 3: 本地生活 - insert a new event row (row 3)
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(3).Insert()

$ws3.Range("A3").Value = 2
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = "2024-07-01"
$ws3.Range("C3").Value = "广州·NIJISANJI EN 官方授权主题店"
$ws3.Range("D3").Value = "天河路299号 时尚天河商业广场"
$ws3.Range("E3").Value = "2024.07.01 00:00-07.15 23:59"
$ws3.Range("F3").Value = 3
$ws3.Range("G3").Value = 30
$ws3.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=86338"
$ws3.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/dB7yQHnF1716795983429.jpeg"

# Match the existing row-index cell style (bold/centered/bordered), like A2
$ws3.Range("A2").Copy()
$ws3.Range("A3").PasteSpecial(-4122) | Out-Null
$ws3.Range("A3").Value = 2

# Re-apply the plain (non-text-forced) formatting used by the rest of the
# column, now that the literal date-like string is safely stored
$ws3.Range("B2").Copy()
$ws3.Range("B3").PasteSpecial(-4122) | Out-Null

# -----------------------------------------------------------------
# Sheet 4: 全部类型 - update "want to go" counts, then insert the
# same new event row (row 28) that was added to 本地生活
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F4").Value  = 1825
$ws4.Range("F5").Value  = 405
$ws4.Range("F6").Value  = 1491
$ws4.Range("F7").Value  = 855
$ws4.Range("F10").Value = 740
$ws4.Range("E11").Value = "2024.06.08 10:00-06.09 16:00"
$ws4.Range("F11").Value = 13187
$ws4.Range("F12").Value = 13063
$ws4.Range("F13").Value = 995
$ws4.Range("F17").Value = 67
$ws4.Range("F18").Value = 639
$ws4.Range("F20").Value = 15
$ws4.Range("F21").Value = 2063
$ws4.Range("F22").Value = 54
$ws4.Range("F26").Value = 18

$ws4.Rows.Item(28).Insert()

$ws4.Range("A28").Value = 27
$ws4.Range("B28").NumberFormat = "@"
$ws4.Range("B28").Value = "2024-07-01"
$ws4.Range("C28").Value = "广州·NIJISANJI EN 官方授权主题店"
$ws4.Range("D28").Value = "天河路299号 时尚天河商业广场"
$ws4.Range("E28").Value = "2024.07.01 00:00-07.15 23:59"
$ws4.Range("F28").Value = 3
$ws4.Range("G28").Value = 30
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86338"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202405/dB7yQHnF1716795983429.jpeg"

$ws4.Range("A27").Copy()
$ws4.Range("A28").PasteSpecial(-4122) | Out-Null
$ws4.Range("A28").Value = 27

# Re-apply the plain (non-text-forced) formatting used by the rest of the
# column, now that the literal date-like string is safely stored
$ws4.Range("B27").Copy()
$ws4.Range("B28").PasteSpecial(-4122) | Out-Null

# Renumber the sequential index column for all rows pushed down by the insert
for ($r = 29; $r -le 36; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# The "Marcin Patrzalek" event row shifted from row 33 to row 34; its
# want-to-go count changes from 0 to 1 (same update as 演出!F9 above)
$ws4.Range("F34").Value = 1
